$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "64.423.64"; E = "  -1.13%  " }
    3  = @{ D = "3.463.24";  E = "  +0.00%  " }
    4  = @{ D = $null;       E = "  +0.29%  " }
    5  = @{ D = "574.22";    E = "  -0.46%  " }
    6  = @{ D = "160.09";    E = "  -0.43%  " }
    7  = @{ D = $null;       E = "  +0.15%  " }
    8  = @{ D = "3.465.06";  E = "  -0.15%  " }
    9  = @{ D = "0.573";     E = "  -5.58%  " }
    10 = @{ D = "7.22";      E = "  -0.66%  " }
    11 = @{ D = $null;       E = "  -2.98%  " }
    12 = @{ D = $null;       E = "  -3.30%  " }
    13 = @{ D = "4.067.32";  E = "  +0.24%  " }
    14 = @{ D = $null;       E = "  -0.10%  " }
    15 = @{ D = "27.53";     E = "  -2.78%  " }
    16 = @{ D = $null;       E = "  -7.76%  " }
    17 = @{ D = "64.661.88"; E = "  -0.72%  " }
    18 = @{ D = "3.459.95";  E = "  +0.69%  " }
    19 = @{ D = "6.21";      E = "  -4.19%  " }
    20 = @{ D = $null;       E = "  -3.42%  " }
    21 = @{ D = "379.72";    E = "  -0.74%  " }
    22 = @{ D = "7.94";      E = "  -2.89%  " }
    23 = @{ D = $null;       E = "  +0.14%  " }
    24 = @{ D = "72.59";     E = "  -0.88%  " }
    25 = @{ D = $null;       E = "  -5.02%  " }
    26 = @{ D = $null;       E = "  -1.21%  " }
    27 = @{ D = "9.82";      E = "  -2.79%  " }
    28 = @{ D = $null;       E = "  +0.20%  " }
    29 = @{ D = $null;       E = "  +0.32%  " }
    30 = @{ D = $null;       E = "  -0.83%  " }
    31 = @{ D = $null;       E = "  -6.64%  " }
    32 = @{ D = $null;       E = "  -1.78%  " }
    33 = @{ D = "23.31";     E = "  -1.66%  " }
    34 = @{ D = $null;       E = "  -3.67%  " }
    35 = @{ D = "1.57";      E = "  -2.97%  " }
    36 = @{ D = "161.30";    E = "  -0.19%  " }
    37 = @{ D = $null;       E = "  -3.60%  " }
    38 = @{ D = "0.824";     E = "  +5.44%  " }
    39 = @{ D = "26.86";     E = "  -0.98%  " }
    40 = @{ D = "0.0745";    E = "  -5.12%  " }
    41 = @{ D = "2.832.73";  E = "  -2.66%  " }
    42 = @{ D = "4.50";      E = "  -5.03%  " }
    43 = @{ D = "42.80";     E = "  -0.83%  " }
    44 = @{ D = "6.44";      E = "  -6.32%  " }
    45 = @{ D = "25.81";     E = "  -0.76%  " }
    46 = @{ D = $null;       E = "  -2.92%  " }
    47 = @{ D = "2.37";      E = "  +7.84%  " }
    48 = @{ D = "333.17";    E = "  +2.97%  " }
    49 = @{ D = $null;       E = "  -3.52%  " }
    50 = @{ D = $null;       E = "  -2.52%  " }
    51 = @{ D = "0.841";     E = "  -4.13%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $vals.D
    }
    $cellE = $ws.Cells.Item($row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $vals.E
}
